# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 50 (pushing the existing rows
# 50-125 down to 51-126) for "Hortaliza, Terminal La Palmera de La Serena -
# Poroto verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 50; this shifts rows
# 50..125 down to 51..126 and grows the sheet dimension to A1:R126.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new observation.
$ws.Range("A50").Value = 8
$ws.Range("B50").Value = "Terminal La Palmera de La Serena"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = 44477
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 100112031
$ws.Range("G50").Value = "Poroto verde"
$ws.Range("H50").Value = "Magnum"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 520
$ws.Range("K50").Value = 37000
$ws.Range("L50").Value = 38000
$ws.Range("M50").Value = 37500
$ws.Range("N50").Value = "$/malla 25 kilos"
$ws.Range("O50").Value = "Perú"
$ws.Range("P50").Value = 1500
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
